$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.921146876822121
$ws.Range("D2").Value = 4.184443087938103
$ws.Range("E2").Value = 11.28924422235743
$ws.Range("F2").Value = 20.02847268923449
$ws.Range("G2").Value = 21.4513911209307
$ws.Range("H2").Value = 12.4387210588638
$ws.Range("K2").Value = 12.88947024689801
$ws.Range("M2").Value = 14.86705293742426
$ws.Range("O2").Value = 17.95422174426258
$ws.Range("C3").Value = 4.749573566289603
$ws.Range("D3").Value = 4.10821793630031
$ws.Range("E3").Value = 11.25935867968669
$ws.Range("F3").Value = 20.09163125322174
$ws.Range("G3").Value = 21.54693890887307
$ws.Range("H3").Value = 12.50400563466551
$ws.Range("K3").Value = 12.19861119684065
$ws.Range("M3").Value = 14.50454833575457
$ws.Range("O3").Value = 18.06137919034512
$ws.Range("C4").Value = 4.642468123253209
$ws.Range("D4").Value = 4.060327656645395
$ws.Range("E4").Value = 11.24675143300196
$ws.Range("F4").Value = 20.13857339558845
$ws.Range("G4").Value = 21.61818501749705
$ws.Range("H4").Value = 12.54695275279085
$ws.Range("K4").Value = 11.75194073444604
$ws.Range("M4").Value = 14.27958291742195
$ws.Range("O4").Value = 18.13318278806957
$ws.Range("C5").Value = 4.598458201387086
$ws.Range("D5").Value = 4.040554756660768
$ws.Range("E5").Value = 11.24305638579636
$ws.Range("F5").Value = 20.15974247355056
$ws.Range("G5").Value = 21.65034865987659
$ws.Range("H5").Value = 12.56517274059843
$ws.Range("K5").Value = 11.56437657168098
$ws.Range("M5").Value = 14.18743915059363
$ws.Range("O5").Value = 18.16394664344131
$ws.Range("C6").Value = 4.591130852377912
$ws.Range("D6").Value = 4.037256435341947
$ws.Range("E6").Value = 11.24252985518763
$ws.Range("F6").Value = 20.16338038768605
$ws.Range("G6").Value = 21.65587746682826
$ws.Range("H6").Value = 12.56824153360977
$ws.Range("K6").Value = 11.53290063544033
$ws.Range("M6").Value = 14.17211440480594
$ws.Range("O6").Value = 18.16914553260619
$ws.Range("C7").Value = 4.641875953508104
$ws.Range("D7").Value = 4.060062011876173
$ws.Range("E7").Value = 11.2466957635785
$ws.Range("F7").Value = 20.13885064784055
$ws.Range("G7").Value = 21.6186061583319
$ws.Range("H7").Value = 12.54719556531192
$ws.Range("K7").Value = 11.74943345611115
$ws.Range("M7").Value = 14.27834195347625
$ws.Range("O7").Value = 18.13359160339798
$ws.Range("C8").Value = 4.862398616130411
$ws.Range("D8").Value = 4.158395977699633
$ws.Range("E8").Value = 11.2777476902326
$ws.Range("F8").Value = 20.04854826246812
$ws.Range("G8").Value = 21.48170600962083
$ws.Range("H8").Value = 12.46063653054127
$ws.Range("K8").Value = 12.65598718748811
$ws.Range("M8").Value = 14.74263362707237
$ws.Range("O8").Value = 17.98991801678979
$ws.Range("C9").Value = 5.277763720310685
$ws.Range("D9").Value = 4.341918405540331
$ws.Range("E9").Value = 11.38414689773185
$ws.Range("F9").Value = 19.93675277740873
$ws.Range("G9").Value = 21.31441541995118
$ws.Range("H9").Value = 12.31365010593656
$ws.Range("K9").Value = 14.25158014306361
$ws.Range("M9").Value = 15.62871573817899
$ws.Range("O9").Value = 17.75618727669348
$ws.Range("C10").Value = 5.568801947326789
$ws.Range("D10").Value = 4.470220568267966
$ws.Range("E10").Value = 11.48980144585261
$ws.Range("F10").Value = 19.89505392338454
$ws.Range("G10").Value = 21.25485924127613
$ws.Range("H10").Value = 12.21959525839431
$ws.Range("K10").Value = 15.30896818871485
$ws.Range("M10").Value = 16.25803116902561
$ws.Range("O10").Value = 17.61419619010441
$ws.Range("C11").Value = 5.697477611171164
$ws.Range("D11").Value = 4.527004778738521
$ws.Range("E11").Value = 11.54372428658055
$ws.Range("F11").Value = 19.88497496107278
$ws.Range("G11").Value = 21.24182081017199
$ws.Range("H11").Value = 12.17984779338031
$ws.Range("K11").Value = 15.7645190269161
$ws.Range("M11").Value = 16.53832155787024
$ws.Range("O11").Value = 17.55615495136856
$ws.Range("C12").Value = 5.74562238840886
$ws.Range("D12").Value = 4.548267359434355
$ws.Range("E12").Value = 11.56497311976842
$ws.Range("F12").Value = 19.88244364447642
$ws.Range("G12").Value = 21.23892424447041
$ws.Range("H12").Value = 12.16523473478668
$ws.Range("K12").Value = 15.93333177662954
$ws.Range("M12").Value = 16.64350204378604
$ws.Range("O12").Value = 17.53512641153859
$ws.Range("C13").Value = 5.735280120216182
$ws.Range("D13").Value = 4.543698970908634
$ws.Range("E13").Value = 11.5603601504369
$ws.Range("F13").Value = 19.88293154567987
$ws.Range("G13").Value = 21.239457026602
$ws.Range("H13").Value = 12.16836239416464
$ws.Range("K13").Value = 15.89713968383784
$ws.Range("M13").Value = 16.62089369450294
$ws.Range("O13").Value = 17.53961288131649
$ws.Range("C14").Value = 5.701450423106205
$ws.Range("D14").Value = 4.52875895291756
$ws.Range("E14").Value = 11.54545590223829
$ws.Range("F14").Value = 19.88474091250562
$ws.Range("G14").Value = 21.24154150559988
$ws.Range("H14").Value = 12.17863677244339
$ws.Range("K14").Value = 15.77848155676841
$ws.Range("M14").Value = 16.54699458235109
$ws.Range("O14").Value = 17.55440582529646
$ws.Range("C15").Value = 5.680651673556384
$ws.Range("D15").Value = 4.51957608573733
$ws.Range("E15").Value = 11.53643420185607
$ws.Range("F15").Value = 19.88601677777746
$ws.Range("G15").Value = 21.24308461339414
$ws.Range("H15").Value = 12.18498726939111
$ws.Range("K15").Value = 15.70531805600966
$ws.Range("M15").Value = 16.50160151849035
$ws.Range("O15").Value = 17.56359095009724
$ws.Range("C16").Value = 5.560313696660612
$ws.Range("D16").Value = 4.466476639950033
$ws.Range("E16").Value = 11.48639430290705
$ws.Range("F16").Value = 19.89589212777232
$ws.Range("G16").Value = 21.25599591050799
$ws.Range("H16").Value = 12.2222541198411
$ws.Range("K16").Value = 15.27868149804437
$ws.Range("M16").Value = 16.23958444907774
$ws.Range("O16").Value = 17.61812186487442
$ws.Range("C17").Value = 5.485502372130129
$ws.Range("D17").Value = 4.433487374701892
$ws.Range("E17").Value = 11.45718868465502
$ws.Range("F17").Value = 19.90423271332833
$ws.Range("G17").Value = 21.26753133435017
$ws.Range("H17").Value = 12.2458954300064
$ws.Range("K17").Value = 15.01040628390424
$ws.Range("M17").Value = 16.07723925039725
$ws.Range("O17").Value = 17.65325870626202
$ws.Range("C18").Value = 5.442125260528529
$ws.Range("D18").Value = 4.414364727357792
$ws.Range("E18").Value = 11.4409429109063
$ws.Range("F18").Value = 19.90986644198566
$ws.Range("G18").Value = 21.27548805239417
$ws.Range("H18").Value = 12.25977919731177
$ws.Range("K18").Value = 14.85370641201152
$ws.Range("M18").Value = 15.98330430346351
$ws.Range("O18").Value = 17.6740846230466
$ws.Range("C19").Value = 5.427380340820759
$ws.Range("D19").Value = 4.407865102232456
$ws.Range("E19").Value = 11.43553763221476
$ws.Range("F19").Value = 19.91191734972039
$ws.Range("G19").Value = 21.27840842848144
$ws.Range("H19").Value = 12.26452906271184
$ws.Range("K19").Value = 14.80024033202237
$ws.Range("M19").Value = 15.95140678375953
$ws.Range("O19").Value = 17.68124148673131
$ws.Range("C20").Value = 5.493502501467937
$ws.Range("D20").Value = 4.437014562402567
$ws.Range("E20").Value = 11.46024057665422
$ws.Range("F20").Value = 19.90325821972947
$ws.Range("G20").Value = 21.266166407851
$ws.Range("H20").Value = 12.24334917068419
$ws.Range("K20").Value = 15.03921278710594
$ws.Range("M20").Value = 16.09457969120828
$ws.Range("O20").Value = 17.64945450107495
$ws.Range("C21").Value = 5.711403165548816
$ws.Range("D21").Value = 4.533153821391836
$ws.Range("E21").Value = 11.54981124814595
$ws.Range("F21").Value = 19.88417452449428
$ws.Range("G21").Value = 21.24087371541059
$ws.Range("H21").Value = 12.1756070270797
$ws.Range("K21").Value = 15.8134348152628
$ws.Range("M21").Value = 16.56872732476323
$ws.Range("O21").Value = 17.55003491868835
$ws.Range("C22").Value = 5.850404076684443
$ws.Range("D22").Value = 4.594579893008486
$ws.Range("E22").Value = 11.61317759899581
$ws.Range("F22").Value = 19.87919655496014
$ws.Range("G22").Value = 21.23624561448777
$ws.Range("H22").Value = 12.13388993076332
$ws.Range("K22").Value = 16.2978897138324
$ws.Range("M22").Value = 16.87297806500218
$ws.Range("O22").Value = 17.49060206609861
$ws.Range("C23").Value = 5.776542803046239
$ws.Range("D23").Value = 4.561928418937288
$ws.Range("E23").Value = 11.57892111601185
$ws.Range("F23").Value = 19.8811657225244
$ws.Range("G23").Value = 21.23762100374553
$ws.Range("H23").Value = 12.15592072105041
$ws.Range("K23").Value = 16.04130717223917
$ws.Range("M23").Value = 16.71113948813147
$ws.Range("O23").Value = 17.52181251416703
$ws.Range("C24").Value = 5.48988678276834
$ws.Range("D24").Value = 4.435420407862408
$ws.Range("E24").Value = 11.45885911790796
$ws.Range("F24").Value = 19.90369617658087
$ws.Range("G24").Value = 21.26677936582917
$ws.Range("H24").Value = 12.24449942485137
$ws.Range("K24").Value = 15.02619704200827
$ws.Range("M24").Value = 16.08674194088127
$ws.Range("O24").Value = 17.65117243513028
$ws.Range("C25").Value = 5.167634243982848
$ws.Range("D25").Value = 4.293358862196595
$ws.Range("E25").Value = 11.35050770853263
$ws.Range("F25").Value = 19.95993651252366
$ws.Range("G25").Value = 21.34865730822322
$ws.Range("H25").Value = 12.31365010593656
$ws.Range("K25").Value = 13.83988604257179
$ws.Range("M25").Value = 15.62871573817899
$ws.Range("O25").Value = 17.81422986130668
